$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Role")

# Add new row 17 data, mirroring row 8 (pointOfContact/contact/RI_414)
# but with a new original_value "businessExpert" in column A.
$ws.Range("A17").Value = "businessExpert"
$ws.Range("B17").Value = $ws.Range("B8").Value2
$ws.Range("C17").Value = $ws.Range("C8").Value2
$ws.Range("D17").Value = $ws.Range("D8").Value2

# Scroll / selection state to match the authored view.
$ws.Range("B17:D17").Select()
$excel.ActiveWindow.ScrollRow = 10
